$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.207.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "'3.591.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("D5").Value = "'581.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").Value = "'192.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").Value = "'0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("D8").Value = "'3.586.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.28%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").Value = "'0.180"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("D11").Value = "'0.664"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "'57.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.33%  "
$ws.Range("D13").Value = "'0.0000305"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("D14").Value = "'9.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "'4.177.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("D16").Value = "'20.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.35%  "
$ws.Range("D17").Value = "'3.605.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.45%  "
$ws.Range("D18").Value = "'70.252.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("D19").Value = "'12.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").Value = "'484.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("D23").Value = "'19.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.18%  "
$ws.Range("E24").Value = "  -10.02%  "
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("D26").Value = "'89.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.48%  "
$ws.Range("D27").Value = "'3.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").Value = "'11.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("E29").Value = "  +1.95%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.42%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'32.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("D32").Value = "'0.122"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.18%  "
$ws.Range("D33").Value = "'12.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("D34").Value = "'66.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("D35").Value = "'608.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").Value = "'40.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.24%  "
$ws.Range("D37").Value = "'0.0₃0809"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.31%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.147"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.404"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  +12.55%  "
$ws.Range("D42").Value = "'3.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "'3.309.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").Value = "'3.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.73%  "
$ws.Range("D45").Value = "'3.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.79%  "
$ws.Range("D46").Value = "'0.0451"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("D47").Value = "'9.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.09%  "
$ws.Range("D48").Value = "'3.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("E51").Value = "  +0.36%  "
